$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.024.71'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '2.948.22'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '379.56'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '101.33'
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").Value = '36.19'
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '18.34'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.406.13'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("E15").Value = '  +4.02%  '
$ws.Range("D16").Value = '12.01'
$ws.Range("E16").Value = '  +69.05%  '
$ws.Range("D17").Value = '2.962.24'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = '0.999'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").Value = '50.972.62'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").Value = '3.06'
$ws.Range("E20").Value = '  -4.22%  '
$ws.Range("D21").Value = '12.50'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").Value = '0.0₃0951'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = '266.70'
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").Value = '3.22'
$ws.Range("E25").Value = '  +12.18%  '
$ws.Range("D26").Value = '8.18'
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").Value = '7.16'
$ws.Range("E27").Value = '  -5.12%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = '25.62'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("D32").Value = '10.09'
$ws.Range("E32").Value = '  +2.93%  '
$ws.Range("D33").Value = '50.51'
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '33.51'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '0.0434'
$ws.Range("E36").Value = '  -5.44%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  +4.36%  '
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '16.62'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '2.55'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("D43").Value = '118.09'
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").Value = '21.55'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").Value = '3.47'
$ws.Range("E45").Value = '  +7.91%  '
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").Value = '2.010.76'
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("E49").Value = '  -4.42%  '
$ws.Range("E50").Value = '  -6.76%  '
$ws.Range("D51").Value = '5.30'
$ws.Range("E51").Value = '  +4.70%  '
